# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
#  - corrige el orden del periodo de mora del trabajador PPT (644798):
#      2307 en la primera fila, 2306 en la segunda -> 2306, 2307
#  - actualiza el Salario Basico de esas dos filas y el de YUNAIDIS
#  - elimina las dos filas de detalle de SHIRLEY EDITH BANQUETH GARCES
#    (el bloque de firma se recorre automaticamente hacia arriba)
#  - refresca los totales de Valor Mora, Cant. Trabajadores y Cant. Periodos

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# La ultima fila de detalle (20) tiene el formato de "cierre de tabla"
# (borde inferior). Como esa fila se va a borrar y la fila 18 (YUNAIDIS)
# pasara a ser la ultima fila de detalle, le copiamos primero ese formato.
$ws.Range("B20:J20").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Corrige el periodo de mora (estaban invertidos) para el trabajador PPT
$ws.Range("E16").Value = "2306"
$ws.Range("E17").Value = "2307"

# Actualiza los valores de Salario Basico
$ws.Range("G16").Value = 1600000
$ws.Range("G17").Value = 1600000
$ws.Range("G18").Value = 1423500

# Elimina las dos filas de detalle de SHIRLEY EDITH BANQUETH GARCES
$ws.Range("B19:J20").EntireRow.Delete()

# Refresca los totales del encabezado
$ws.Range("E11").Value = 120533
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3

# El ancho de la columna D (Nombre Trabajador) se reajusta porque el
# texto mas largo que contenia ("SHIRLEY EDITH BANQUETH GARCES") fue
# eliminado.
$ws.Columns("D").ColumnWidth = 30.5
